# Swap "X, System" -> "System, X" in column G (Recorded By) for every row
# on the active sheet. This reorders the comma-separated list of recorders
# so that "System" always appears first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*, System") {
        $prefix = $val.Substring(0, $val.Length - 8)  # strip ", System"
        $cell.Value2 = "System, " + $prefix
    }
}
